# Regenerate save_data to use K (strikeouts) instead of Strike# in column G.
# This updates the "K" column (column G, rows 2-73) with the recalculated
# strikeout counts for familia_jeurys 2021 save data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K values for rows 2..73 (one entry per row, in order)
$newK = @(2,1,0,1,1,0,1,1,0,0,0,1,2,1,2,1,3,1,0,2,2,0,1,0,0,0,3,1,1,1,1,3,2,3,1,2,0,0,2,3,2,0,2,0,0,1,2,1,1,1,1,2,1,2,0,1,1,1,2,1,0,0,2,1,1,0,1,0,0,1,1,0)

$startRow = 2
for ($i = 0; $i -lt $newK.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 7).Value = $newK[$i]
}
